# KRST-82, Feature, Export Student List Info
# Applies the "学生信息批量导出" workbook edits:
#  - rename the helper sheet and add a checksum/date-format block to it
#  - tidy up sheet1's title row (height/style) and view (drop frozen panes)
#  - reorder the data-validation rules
#  - make the helper sheet the active tab

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 学生信息
$ws2 = $wb.Worksheets.Item(2)   # was "Sheet1"

# ---------------------------------------------------------------------------
# Sheet1 ("学生信息"): title row height + view cleanup
# ---------------------------------------------------------------------------

# Title row is shorter in the new template (26.4 -> 18) and no longer needs an
# explicit custom-height flag beyond the plain height set below.
$ws1.Rows.Item(1).RowHeight = 18

# Drop the frozen header panes / multi-pane selection in favour of a single
# plain selection at G3, and make sure sheet1 is not left as the tab-selected
# sheet (sheet2 becomes the active one further down).
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("G3").Select()

# Data validations: re-created in the new order (text-length rule first, then
# the M/N-list rule second) and the M/N sqref is merged into one range.
$ws1.Cells.Validation.Delete()
$ws1.Range("M3:N1048576").Validation.Add(6, 1, 1, 8, 11)
$ws1.Range("D3:D1048576").Validation.Add(3, 1, 1, '"男,女"')

# ---------------------------------------------------------------------------
# Sheet2: rename to the "helper info" sheet and populate the checksum block
# ---------------------------------------------------------------------------

$ws2.Name = "辅助信息（请勿操作）"

$ws2.Range("A1").Value = "校验码"
$ws2.Range("A2").Value = "H86D`$8#a"
$ws2.Range("B1").Value = "日期格式"
$ws2.Range("B2").Value = 44562
$ws2.Range("B2").NumberFormat = "yyyy\-mm\-dd;@"

$ws2.Columns.Item(2).ColumnWidth = 9.5

$ws2.Range("C2").Select()

# Make the helper sheet the active/visible tab, matching the exported
# template default.
$ws2.Activate()
